# Regenerate save_data: column G ("K") is recalculated (std/mean, s_vals)
# and replaces the previous "Strike#" values. Write the freshly computed
# K values for each data row back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    22 = 3
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 2
    35 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    44 = 0
    45 = 2
    46 = 0
    48 = 1
    49 = 2
    50 = 1
    51 = 2
    52 = 0
    53 = 1
    54 = 0
    55 = 3
    56 = 1
    58 = 2
    59 = 1
    60 = 2
    61 = 1
    62 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
